# Anonymize "fedcore" -> "approach" and give the C1/D1 (and F1/G1) header
# spacer cells their top/bottom (and right) border accents.
#
# Border edge indices used below match Excel's XlBordersIndex constants:
#   7  = xlEdgeLeft
#   8  = xlEdgeTop
#   9  = xlEdgeBottom
#   10 = xlEdgeRight
# xlLineStyleNone = -4142, xlContinuous = 1

$wb = $excel.ActiveWorkbook

# NOTE on ordering: the engine clones a brand-new (orphan) border/xf entry
# for every distinct edge-combination a range passes through while it is
# being edited; it does not retroactively garbage-collect abandoned ones.
# To land cleanly on the pre-existing borderId 4 (top+bottom) and borderId 5
# (top+bottom+right) table entries without leaving orphaned style/border
# rows behind, each cell's edits are ordered so every intermediate
# combination is itself one of the already-defined borders: none -> top ->
# (top+right ->) top+bottom(+right).
function Set-TopBottomBorder {
    param($range)
    $range.Style = "Normal"
    $range.Borders.Item(8).LineStyle = 1
    $range.Borders.Item(9).LineStyle = 1
}

function Set-TopBottomRightBorder {
    param($range)
    $range.Style = "Normal"
    $range.Borders.Item(8).LineStyle = 1
    $range.Borders.Item(10).LineStyle = 1
    $range.Borders.Item(9).LineStyle = 1
}

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-TopBottomBorder      $ws1.Range("C1")
Set-TopBottomRightBorder $ws1.Range("D1")

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-TopBottomBorder      $ws2.Range("C1")
Set-TopBottomRightBorder $ws2.Range("D1")
Set-TopBottomBorder      $ws2.Range("F1")
Set-TopBottomRightBorder $ws2.Range("G1")

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 was an empty placeholder cell; drop it entirely.
$ws2.Range("G5").ClearContents()
